$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-20
$data = @(
    @(1, 7),
    @(1, 5),
    @(5, 6),
    @(1, 5),
    @(1, 6),
    @(1, 7),
    @(1, 6),
    @(6, 8),
    @(8, 9),
    @(5, 5),
    @(5, 7),
    @(6, 8),
    @(7, 9),
    @(10, 10),
    @(6, 7),
    @(1, 5),
    @(1, 4),
    @(1, 3),
    @(4, 6)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
